# Apply updated cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.419.07"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "1.793.53"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'339.30"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").Value = "'0.9971"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").Value = "'0.3920"
$ws.Range("E7").Value = "  +2.52%  "
$ws.Range("D8").Value = "'0.3467"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("D9").Value = "'48.33"
$ws.Range("E9").Value = "  -2.50%  "
$ws.Range("D10").Value = "1.194"
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("D11").Value = "'0.07482"
$ws.Range("E11").Value = "  -2.78%  "
$ws.Range("D12").Value = "'0.9960"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").Value = "'21.87"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D15").Value = "1.795.75"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("D16").Value = "'7.147"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "'0.00001098"
$ws.Range("E17").Value = "  -1.77%  "
$ws.Range("D19").Value = "'84.76"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("D20").Value = "'0.9960"
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("D21").Value = "'17.70"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").Value = "'6.563"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").Value = "27.477.77"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").Value = "12.42"
$ws.Range("E24").Value = "  -4.80%  "
$ws.Range("D25").Value = "2.399"
$ws.Range("E25").Value = "  -2.79%  "
$ws.Range("D28").Value = "1.466"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "'157.50"
$ws.Range("E29").Value = "  +4.04%  "
$ws.Range("D30").Value = "2.000.13"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").Value = "'135.68"
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("D32").Value = "'4.018"
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("D33").Value = "'6.047"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("D34").Value = "'0.08748"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").Value = "13.08"
$ws.Range("E35").Value = "  -4.90%  "
$ws.Range("D36").Value = "'1.624"
$ws.Range("E36").Value = "  -3.76%  "
$ws.Range("D40").Value = "'0.6846"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("D41").Value = "'0.2216"
$ws.Range("E41").Value = "  -1.51%  "
$ws.Range("D42").Value = "'1.255"
$ws.Range("E42").Value = "  -3.29%  "
$ws.Range("D43").Value = "'8.412"
$ws.Range("E43").Value = "  -7.12%  "
$ws.Range("D44").Value = "'14.61"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D47").Value = "'3.870"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").Value = "'2.140"
$ws.Range("E48").Value = "  -1.35%  "
$ws.Range("D49").Value = "'131.99"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").Value = "'0.07187"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D18").Value = "'0.06678"
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.520"
$ws.Range("E26").Value = "  -4.57%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'21.29"
$ws.Range("E27").Value = "  -3.31%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'5.450"
$ws.Range("E37").Value = "  -2.51%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02421"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06501"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6395"
$ws.Range("E45").Value = "  -2.15%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "'0.9947"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'79.51"
$ws.Range("E51").Value = "  -1.42%  "
